$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Byars, 95% CI, x = 1): amended formula results
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 0.013072132901504351
$ws.Range("D4").Value = 5.5637555824752605

# Row 7 (Byars, 99.8% CI, x = 10): amended formula results
$ws.Range("A7").Value = 10
$ws.Range("C7").Value = 2.916318374695686
$ws.Range("D7").Value = 24.193724927880535

# Row 14 (new Byars, 99.8% CI test case, x = 1, lowercl blanked out as negative byars result)
$ws.Range("A14").Value = 1
$ws.Range("D14").Value = 9.362180217908362
$ws.Range("E14").Value = "99.8%"
$ws.Range("F14").Value = "Byars"

# Update the active selection to match the authored workbook state
$ws.Range("D21").Select() | Out-Null
